# "Generate Report for Handback"
#
# The localization-status report is refreshed after a handback completes
# for f607f7c6-f204-4dde-82d2-c1cadd15fb77 (row 2 on every sheet):
#   - Overview: status columns (zh-cn / de-de) flip from "Ready for
#     handoff" to "Handed back: in sync with en-US".
#   - Each language sheet (zh-cn, de-de) records the handed-back target
#     file (with a hyperlink, like the Source File Name column), the
#     handback xliff file name, and the handback datetime.
#   - A few columns are widened so the new, longer values are legible.

$wb = $excel.ActiveWorkbook

$targetMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/de928e5ac20f8f0d30de3bacde0dae1bb40af6ad/e2e/f607f7c6-f204-4dde-82d2-c1cadd15fb77.md"
$hyperlinkColor = 15570276  # decimal BGR for RGB FF6495ED, matching the workbook's existing HyperLink style

# ---------------------------------------------------------------------
# Overview sheet: row 2 (f607f7c6) is now handed back & in sync for
# both locales.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"

$overview.Columns.Item(5).ColumnWidth = 29.2
$overview.Columns.Item(6).ColumnWidth = 29.2

# ---------------------------------------------------------------------
# Per-locale sheets: fill in Latest Target File / Latest Handback File /
# Latest Handback DateTime for the f607f7c6 row, and widen the Status /
# Latest Target File / Latest Handback File columns.
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("I2").Value = "f607f7c6-f204-4dde-82d2-c1cadd15fb77.md"
$zhcn.Range("J2").Value = "f607f7c6-f204-4dde-82d2-c1cadd15fb77.6fe5a1463f5ee4c99bd7d3af2e2db969f44e613b.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-23 14:47:41"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $targetMdUrl, $null, $null, "f607f7c6-f204-4dde-82d2-c1cadd15fb77.md")
$zhcn.Range("I2").Font.Underline = 2
$zhcn.Range("I2").Font.Color = $hyperlinkColor

$zhcn.Columns.Item(3).ColumnWidth = 29.2
$zhcn.Columns.Item(9).ColumnWidth = 39.2
$zhcn.Columns.Item(10).ColumnWidth = 39.2

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("I2").Value = "f607f7c6-f204-4dde-82d2-c1cadd15fb77.md"
$dede.Range("J2").Value = "f607f7c6-f204-4dde-82d2-c1cadd15fb77.6fe5a1463f5ee4c99bd7d3af2e2db969f44e613b.de-de.xlf"
$dede.Range("K2").Value = "2016-08-23 14:47:49"

$dede.Hyperlinks.Add($dede.Range("I2"), $targetMdUrl, $null, $null, "f607f7c6-f204-4dde-82d2-c1cadd15fb77.md")
$dede.Range("I2").Font.Underline = 2
$dede.Range("I2").Font.Color = $hyperlinkColor

$dede.Columns.Item(3).ColumnWidth = 29.2
$dede.Columns.Item(9).ColumnWidth = 39.2
$dede.Columns.Item(10).ColumnWidth = 39.2
